{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Change 1: In the \"3\u0410. \u041f\u0440\u043e\u0445\u043e\u0434 \u043f\u043e \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u0435.\" alternative flow, a new\n// bullet step is inserted at the top: \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\u0438\u0432\u0430\u0435\u0442, \u0447\u0442\u043e \u0442\u0438\u043f \u043a\u0430\u0440\u0442\u044b \u2013\n// \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u0430\u044f \u043a\u0430\u0440\u0442\u0430.\" (typed/split the way Word captures an in-place edit,\n// leaving a collapsed \"_GoBack\" bookmark at the point the author paused\n// while typing). The three steps that used to follow it shift down by one,\n// and the step that used to be last (\"\u0423\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u0435 \u043f\u0435\u0440\u0435\u0434\u0430\u0451\u0442\u0441\u044f \u043d\u0430 \u0448\u0430\u0433 4\n// \u043e\u0441\u043d\u043e\u0432\u043d\u043e\u0433\u043e \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u044f.\") ends up as a brand new paragraph with the same\n// bullet formatting as the others (unchanged text/content).\n//\n// Change 2: Further down, in the postconditions of \"\u0421\u0438\u043d\u0445\u0440\u043e\u043d\u0438\u0437\u0438\u0440\u043e\u0432\u0430\u0442\u044c \u0432\u0440\u0435\u043c\u044f \u0438\n// \u0434\u0430\u0442\u0443\", the leftover collapsed \"_GoBack\" bookmark that used to split\n// \"...\u043d\u0430 \u0441\u0442\u0430\u043d\u0446\u0438\u044e\" and \".\" into two runs is removed, leaving one run with the\n// full sentence.\n//\n// NOTE: Change 2 is applied BEFORE change 1 so the pre-existing \"_GoBack\"\n// bookmark is cleared out of the way before a new one is created further up\n// in the document (bookmark names are effectively unique, like in real\n// Word, so clearing the old one first avoids any ambiguity about which\n// \"_GoBack\" a later lookup/deletion would touch).\n\nconst SENTENCE =\n  \"\u044e\u0431\u043e\u043c \u0441\u043b\u0443\u0447\u0430\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u0442\u0441\u044f, \u0447\u0442\u043e \u0434\u0430\u0442\u0430 \u0438 \u0432\u0440\u0435\u043c\u044f \u043d\u0435 \u0431\u0443\u0434\u0443\u0442 \u0441\u0438\u043d\u0445\u0440\u043e\u043d\u0438\u0437\u0438\u0440\u043e\u0432\u0430\u043d\u044b \u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0445\u043e\u0434\u0430 \u043f\u0430\u0441\u0441\u0430\u0436\u0438\u0440\u0430 \u043d\u0430 \u0441\u0442\u0430\u043d\u0446\u0438\u044e.\";\n\nconst sentenceHits = context.document.body.search(SENTENCE, { matchCase: true });\nsentenceHits.load(\"items\");\nawait context.sync();\n\nif (sentenceHits.items.length > 0) {\n  // Re-insert the same text so the run that used to be split by the\n  // bookmark is rebuilt as a single contiguous run.\n  sentenceHits.items[0].insertText(SENTENCE, Word.InsertLocation.replace);\n  await context.sync();\n}\n\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // No \"_GoBack\" bookmark present \u2014 nothing to remove.\n}\n\n// --- Change 1: insert the new step before the existing first step -------\nconst OLD_FIRST_STEP =\n  \"\u0417\u0430\u043f\u0443\u0441\u043a\u0430\u0435\u0442\u0441\u044f \u043e\u0441\u043d\u043e\u0432\u043d\u043e\u0439 \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u0439 \u0432\u043a\u043b\u044e\u0447\u0451\u043d\u043d\u043e\u0433\u043e \u0412\u0418 \u00ab\u041e\u043f\u043b\u0430\u0442\u0438\u0442\u044c \u043f\u0440\u043e\u0445\u043e\u0434 \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u043e\u0439\u00bb.\";\nconst NEW_FIRST_STEP =\n  \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\u0438\u0432\u0430\u0435\u0442, \u0447\u0442\u043e \u0442\u0438\u043f \u043a\u0430\u0440\u0442\u044b \u2013 \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u0430\u044f \u043a\u0430\u0440\u0442\u0430.\";\nconst SPLIT_AFTER = \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\"; // where the author's cursor paused\n\nconst hits = context.document.body.search(OLD_FIRST_STEP, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Target paragraph not found: \" + OLD_FIRST_STEP);\n}\n\nconst target = hits.items[0].paragraphs.getFirst();\ntarget.insertParagraph(NEW_FIRST_STEP, Word.InsertLocation.before);\nawait context.sync();\n\n// Re-find the freshly inserted text and drop a collapsed bookmark right\n// after \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\", matching where Word leaves \"_GoBack\" when the\n// document is saved mid-edit. The new paragraph is the first occurrence of\n// this prefix in the body (the others, further down, already existed).\nconst splitHits = context.document.body.search(SPLIT_AFTER, { matchCase: true });\nsplitHits.load(\"items\");\nawait context.sync();\n\nif (splitHits.items.length > 0) {\n  const splitPoint = splitHits.items[0].getRange(Word.RangeLocation.end);\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Change 1: In the \"3\u0410. \u041f\u0440\u043e\u0445\u043e\u0434 \u043f\u043e \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u0435.\" alternative flow, a new\n# bullet step is inserted at the top: \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\u0438\u0432\u0430\u0435\u0442, \u0447\u0442\u043e \u0442\u0438\u043f \u043a\u0430\u0440\u0442\u044b \u2013\n# \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u0430\u044f \u043a\u0430\u0440\u0442\u0430.\" (typed/split the way Word captures an in-place edit,\n# leaving a collapsed \"_GoBack\" bookmark at the point the author paused while\n# typing). The three steps that used to follow it shift down by one, and the\n# step that used to be last (\"\u0423\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u0435 \u043f\u0435\u0440\u0435\u0434\u0430\u0451\u0442\u0441\u044f \u043d\u0430 \u0448\u0430\u0433 4 \u043e\u0441\u043d\u043e\u0432\u043d\u043e\u0433\u043e\n# \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u044f.\") ends up as a brand new paragraph with the same bullet\n# formatting as the others (unchanged text/content).\n#\n# Change 2: Further down, in the postconditions of \"\u0421\u0438\u043d\u0445\u0440\u043e\u043d\u0438\u0437\u0438\u0440\u043e\u0432\u0430\u0442\u044c \u0432\u0440\u0435\u043c\u044f \u0438\n# \u0434\u0430\u0442\u0443\", the leftover collapsed \"_GoBack\" bookmark that used to split\n# \"...\u043d\u0430 \u0441\u0442\u0430\u043d\u0446\u0438\u044e\" and \".\" into two runs is removed, leaving one run with the\n# full sentence.\n#\n# NOTE: Change 2 is applied BEFORE change 1 so the pre-existing \"_GoBack\"\n# bookmark is cleared out of the way before a new one is created further up\n# in the document (bookmark names must be unique, so clearing the old one\n# first avoids any ambiguity about which \"_GoBack\" a later lookup would\n# touch).\n\n$d = $word.ActiveDocument\n\n# --- Change 2: drop the stale \"_GoBack\" bookmark and rebuild the sentence\n# it used to split into a single run. --------------------------------------\n$sentence = \"\u044e\u0431\u043e\u043c \u0441\u043b\u0443\u0447\u0430\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u0442\u0441\u044f, \u0447\u0442\u043e \u0434\u0430\u0442\u0430 \u0438 \u0432\u0440\u0435\u043c\u044f \u043d\u0435 \u0431\u0443\u0434\u0443\u0442 \u0441\u0438\u043d\u0445\u0440\u043e\u043d\u0438\u0437\u0438\u0440\u043e\u0432\u0430\u043d\u044b \u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0445\u043e\u0434\u0430 \u043f\u0430\u0441\u0441\u0430\u0436\u0438\u0440\u0430 \u043d\u0430 \u0441\u0442\u0430\u043d\u0446\u0438\u044e.\"\n\n$find = $d.Content.Find\n$find.Text = $sentence\n$find.Replacement.Text = $sentence\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\ntry {\n    $d.Bookmarks(\"_GoBack\").Delete()\n} catch {\n    # No \"_GoBack\" bookmark present - nothing to remove.\n}\n\n# --- Change 1: insert the new first step before the existing one ---------\n$oldFirstStep = \"\u0417\u0430\u043f\u0443\u0441\u043a\u0430\u0435\u0442\u0441\u044f \u043e\u0441\u043d\u043e\u0432\u043d\u043e\u0439 \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u0439 \u0432\u043a\u043b\u044e\u0447\u0451\u043d\u043d\u043e\u0433\u043e \u0412\u0418 \u00ab\u041e\u043f\u043b\u0430\u0442\u0438\u0442\u044c \u043f\u0440\u043e\u0445\u043e\u0434 \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u043e\u0439\u00bb.\"\n$newFirstStep = \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\u0438\u0432\u0430\u0435\u0442, \u0447\u0442\u043e \u0442\u0438\u043f \u043a\u0430\u0440\u0442\u044b \u2013 \u0431\u0430\u043d\u043a\u043e\u0432\u0441\u043a\u0430\u044f \u043a\u0430\u0440\u0442\u0430.\"\n$splitAfter = \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\"\n\n$find2 = $d.Content.Find\n$find2.Text = $oldFirstStep\n$found = $find2.Execute()\n\nif ($found) {\n    $targetRange = $find2.Parent\n    # Appending a carriage return inserts a brand new paragraph before the\n    # matched one, inheriting its paragraph formatting (style/numbering),\n    # and sets its text in a single step.\n    $targetRange.InsertBefore(\"$newFirstStep`r\")\n}\n\n# Drop a collapsed \"_GoBack\" bookmark right after \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\" in the\n# newly-inserted paragraph (the first occurrence of that prefix in the\n# document body, since the paragraph we just added precedes every other\n# \"\u0422\u0443\u0440\u043d\u0438\u043a\u0435\u0442 \u043e\u0431\u043d\u0430\u0440\u0443\u0436\u0438\u0432\u0430\u0435\u0442...\" bullet).\n$find3 = $d.Content.Find\n$find3.Text = $splitAfter\n$found3 = $find3.Execute()\nif ($found3) {\n    $splitPoint = $d.Range($find3.Parent.End, $find3.Parent.End)\n    $d.Bookmarks.Add(\"_GoBack\", $splitPoint)\n}\n"}
